$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.213.50'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.895.62'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '''245.75'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("E6").Value = '  +8.16%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = '''40.76'
$ws.Range("E8").Value = '  -3.41%  '
$ws.Range("D9").Value = '''0.347'
$ws.Range("E9").Value = '  +2.93%  '
$ws.Range("D10").Value = '''52.24'
$ws.Range("E10").Value = '  +8.79%  '
$ws.Range("D11").Value = '''0.0719'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").Value = '2.170.44'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").Value = '''12.78'
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = '''0.705'
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").Value = '1.889.93'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Value = '''4.79'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '35.194.91'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").Value = '''71.93'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '0.0₃0817'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '''240.43'
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").Value = '''12.65'
$ws.Range("E22").Value = '  +1.85%  '
$ws.Range("D23").Value = '''4.79'
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D25").Value = '''2.33'
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("D27").Value = '''167.64'
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").Value = '''8.60'
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").Value = '''18.31'
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("D30").Value = '''0.130'
$ws.Range("E30").Value = '  +3.83%  '
$ws.Range("E31").Value = '  +20.01%  '
$ws.Range("D32").Value = '''4.15'
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("D33").Value = '''0.0567'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '''1.90'
$ws.Range("E34").Value = '  +9.04%  '
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '''4.11'
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.917'
$ws.Range("E37").Value = '  -6.27%  '
$ws.Range("E38").Value = '  +13.48%  '
$ws.Range("D39").Value = '''2.02'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''16.42'
$ws.Range("E40").Value = '  +6.19%  '
$ws.Range("E41").Value = '  -1.24%  '
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").Value = '''0.0636'
$ws.Range("E43").Value = '  +7.61%  '
$ws.Range("D44").Value = '''90.00'
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").Value = '1.346.63'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("D46").Value = '''2.42'
$ws.Range("E46").Value = '  +3.15%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '''2.78'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").Value = '''45.64'
$ws.Range("E49").Value = '  -11.65%  '
$ws.Range("D50").Value = '''12.02'
$ws.Range("E50").Value = '  -5.52%  '
$ws.Range("D51").Value = '''6.46'
$ws.Range("E51").Value = '  -2.86%  '
